$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$shp = $m.Shapes.Item(3)
$tr2 = $shp.TextFrame2.TextRange
Write-Output $tr2.Text
$tr2 | Get-Member | Out-String | Write-Output
